# penbmi for template 1
# - cohencu (row 9) now has a result for "Our" instead of "Failed"
# - the "inv=" prefix is dropped from the cav13-2 template text (B8)
# - the (always empty) invSDP column is dropped from the Template 1 table
#
# penbmi for template 1
# - the PresetN column and the trailing "4次以下无解" notes column are
#   dropped from the Template 2 table, shifting Z3 / PolySynth / invSDP
#   one column to the left

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Template 1 table (rows 1-12) ---------------------------------------

# cav13-2 formula text: drop the "inv=" prefix
$ws.Range("B8").Value = "a0+a1*x1+a2*x2`na3*x1^2+a4*x2^2+a5*x1*x2"

# cohencu "Our" result is now available
$ws.Range("C9").Value = 1.33

# Drop the (empty) invSDP column for this table; row 12 (outside the
# table) keeps its own G12 content untouched.
$ws.Range("G1:G11").Clear()

# --- Template 2 table (rows 16-26) ---------------------------------------

# Drop the PresetN column (C) by shifting D:G left into C:F, preserving
# each cell's own formatting; then clear the now-stale trailing columns
# (the old invSDP column G and the notes column H).
$ws.Range("D16:G26").Copy($ws.Range("C16:F26"))
$ws.Range("G16:H26").Clear()

# --- View bookkeeping -----------------------------------------------------
$ws.Application.Goto($ws.Range("A8"), $false)
$ws.Range("F26").Select()
